$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 11 with new staff record, replacing MOUNICA A / VEC-003-04-600 data
$ws.Range("B11").Value = "Lab Instructor"
$ws.Range("J11").Value = "VEC-003-05-3"
$ws.Range("A11").Value = "KUMAR V"
$ws.Range("C11").Value = "/static/images/profile_photos/003/VEC-003-05-3.webp"

# Apply an all-around thin border to the Name/Designation cells (A11:B11)
$ws.Range("A11:B11").Borders.LineStyle = 1
$ws.Range("A11:B11").Borders.Weight = 2

# Update selection to K15 to match final state
$ws.Range("K15").Select()
